$wb = $excel.ActiveWorkbook

# Mapping of cell (column F, "想去人数") updates that apply identically
# to both the "展览" and "全部类型" worksheets.
$updates = @{
    "F3"  = 4987
    "F5"  = 7254
    "F8"  = 96
    "F9"  = 593
    "F12" = 4246
    "F13" = 1704
    "F14" = 95
    "F15" = 91
    "F16" = 2843
    "F19" = 194
    "F20" = 445
    "F21" = 411
    "F22" = 436
    "F23" = 270
    "F24" = 78
    "F25" = 1665
    "F26" = 1136
    "F28" = 1329
    "F32" = 508
    "F33" = 19
    "F34" = 47
    "F35" = 101
    "F36" = 2641
    "F37" = 676
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
